$wb = $excel.ActiveWorkbook

# Rename the worksheet "Sheet1" to "Research"
$wsResearch = $wb.Worksheets.Item("Sheet1")
$wsResearch.Name = "Research"

# Make the "Data" sheet the active/selected sheet (was "Research"/"Sheet1" before),
# and move the selection on that sheet from Q26 to Q4.
$wsData = $wb.Worksheets.Item("Data")
$wsData.Activate() | Out-Null
$wsData.Range("Q4").Select() | Out-Null
